# issue #5: add legislator_id, name, date into dataframe
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$legislatorDate = "2013-12-12"
$legislatorName = "張慶忠"
$legislatorId = 1347

# ----- Header row (row 1): date / legislator_name / legislator_id -----
$h1 = $ws.Cells.Item(1, 8)
$h1.Value = "date"
$h1.Font.Bold = $true
$h1.HorizontalAlignment = -4108
$h1.VerticalAlignment = -4160
$h1.Borders.LineStyle = 1

$i1 = $ws.Cells.Item(1, 9)
$i1.Value = "legislator_name"
$i1.Font.Bold = $true
$i1.HorizontalAlignment = -4108
$i1.VerticalAlignment = -4160
$i1.Borders.LineStyle = 1

$j1 = $ws.Cells.Item(1, 10)
$j1.Value = "legislator_id"
$j1.Font.Bold = $true
$j1.HorizontalAlignment = -4108
$j1.VerticalAlignment = -4160
$j1.Borders.LineStyle = 1

# Find last used data row (column A has the row id, data starts at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    # date column - force text so "2013-12-12" isn't auto-parsed into a date serial
    $hc = $ws.Cells.Item($r, 8)
    $hc.NumberFormat = "@"
    $hc.Value = $legislatorDate
    $hc.Style = "Normal"

    $ic = $ws.Cells.Item($r, 9)
    $ic.Value = $legislatorName
    $ic.Style = "Normal"

    $jc = $ws.Cells.Item($r, 10)
    $jc.Value = $legislatorId
    $jc.Style = "Normal"
}
